$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (B2 changes from 2 to 3; A2 and C2 stay the same)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2125000000

# Add new supply rows 3-7
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1955000000

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 2975000000

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 3570000000

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1360000000

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 1530000000

# Move the selection/active cell to C6, matching the target workbook state
$ws.Range("C6").Select()
